$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: "active" header, "true" (text) for every existing data row.
$ws.Cells.Item(1, 5).Value = "active"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 5)
    # Leading apostrophe forces Excel to store this as literal text ("true")
    # instead of auto-converting to the Boolean TRUE.
    $c.Value = "'true"
    # Drop the quote-prefix formatting flag that the apostrophe trick leaves
    # behind, so the new cells keep the workbook's default (no) style.
    $c.ClearFormats()
}

Write-Output "lastRow=$lastRow"
